# Powerpoint writer: consolidate text run nodes.
# Merge each word with its trailing space into a single run, reducing
# the number of <a:r> nodes emitted for these two text bodies.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title 1: "A" + " " + "slide"  ->  "A " + "slide"
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Characters(1, 2).Text = "A "

# TextBox 3: "Followed" + " " + "by" + " " + "a" + " " + "picture"
#            -> "Followed " + "by " + "a " + "picture"
$caption = $s.Shapes.Item(4).TextFrame.TextRange
$caption.Characters(1, 9).Text = "Followed "
$caption.Characters(10, 3).Text = "by "
$caption.Characters(13, 2).Text = "a "
